$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark "Balance this thang" (row 6) as done: strikethrough B6:D6, keep existing alignment ---
$ws.Range("B6:D6").Font.Strikethrough = $true

# --- New feature row 8: "Debug this bitch" ---
$ws.Range("B8").Value = "Debug this bitch"
$ws.Range("C8").Value = "Garrett"
$d8text = "player transition to offensive positions, players jump into place`nplayer transitioning to defensive positions happens after ai transitions to offence, should be same time`nplayer setter can still be moved around during offensive phase`nwhen calculating the closest pawn, use the updated location and not the indicators location`nstupid ass bug where sometimes after a while of playing, the ball does not get set properly to the servers location"
$ws.Range("D8").Value = $d8text
$ws.Range("B8:C8").VerticalAlignment = -4108
$ws.Range("D8").WrapText = $true
# Colour the two "done/observation" lines green, leave the rest default (automatic/black)
$ws.Range("D8").Characters(1, 65).Font.Color = 5296274
$ws.Range("D8").Characters(67, 105).Font.Color = 5296274
$ws.Rows.Item(8).RowHeight = 72

# --- New feature row 9: "Setter updates" ---
$ws.Range("B9").Value = "Setter updates"
$ws.Range("C9").Value = "Garrett"
$ws.Range("D9").Value = "when the setter makes a dig, have someone else come in and set the ball"
$ws.Range("B9:C9").VerticalAlignment = -4108
$ws.Range("D9").WrapText = $true

# --- New feature row 10: "Serve Receive Rotation" ---
$ws.Range("B10").Value = "Serve Receive Rotation"
$ws.Range("C10").Value = "Garrett"
$d10text = "Make it so that whenever the player is changing around their serve receive positions, they are unable to place their players out of rotation.`nLikely need a function to check whether or not an attempted placement is valid, and then not allowing the player to move beyond that space, similar to the limited movement during reaction phases`nneed some sort of prompt to the player that they attempted to place their piece out of rotation. "
$ws.Range("D10").Value = $d10text
$ws.Range("B10:C10").VerticalAlignment = -4108
$ws.Range("D10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 57.6

# --- Update selection to match where the author ended up working ---
$ws.Range("D8").Select()
